$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Re-ran training (encoder frozen) -> refresh the per-epoch accuracy readings in column B
$ws.Range("B2").Value = 0.9375
$ws.Range("B4").Value = 0.859375
$ws.Range("B5").Value = 0.8125
$ws.Range("B8").Value = 0.765625
$ws.Range("B9").Value = 0.703125
$ws.Range("B10:B11").Value = 0.734375
$ws.Range("B12").Value = 0.6875
$ws.Range("B13").Value = 0.640625
$ws.Range("B14").Value = 0.609375
$ws.Range("B15").Value = 0.703125
$ws.Range("B17").Value = 0.640625
$ws.Range("B18").Value = 0.65625
$ws.Range("B19:B20").Value = 0.609375
$ws.Range("B22").Value = 0.609375
$ws.Range("B23").Value = 0.515625
$ws.Range("B24").Value = 0.5625
$ws.Range("B25").Value = 0.53125
$ws.Range("B26:B27").Value = 0.546875
$ws.Range("B28:B29").Value = 0.53125
$ws.Range("B30").Value = 0.515625
$ws.Range("B31").Value = 0.53125
$ws.Range("B32:B45").Value = 0.546875
$ws.Range("B46:B55").Value = 0.53125
$ws.Range("B56:B69").Value = 0.515625
$ws.Range("B70:B88").Value = 0.5
$ws.Range("B89:B102").Value = 0.484375
$ws.Range("B103").Value = 0.5
$ws.Range("B104").Value = 0.4375
$ws.Range("B105").Value = 0.46875
$ws.Range("B106").Value = 0.421875
$ws.Range("B107").Value = 0.484375
$ws.Range("B108").Value = 0.375
$ws.Range("B109").Value = 0.546875
$ws.Range("B110").Value = 0.46875
$ws.Range("B111").Value = 0.40625
$ws.Range("B112").Value = 0.421875
$ws.Range("B114").Value = 0.546875
$ws.Range("B115").Value = 0.46875
$ws.Range("B116").Value = 0.515625
$ws.Range("B117").Value = 0.53125

# The DisplayOutputs repr cells in column A were regenerated with a new Python object id
$ws.Range("A102:A118").Value = "<__main__.DisplayOutputs object at 0x7f47a0583910>"
